$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetingen frontend & backend")

# --- Insert a new "content length" row above the old row 4 -----------------
# This pushes every row from the old 4..34 down by one (formulas auto-adjust).
$ws.Rows.Item(4).Insert()

# The former "1 level" row under "mount tot data inladen x sneller dan zonder
# cache" (old row 24) is now at row 25 and is a duplicate that must go away.
$ws.Rows.Item(25).Delete()

# New row 4 content: content length measurement
$ws.Range("A4").Value = "content length"
$ws.Range("B4").Value = 4302

# --- Updated measurements: Top Level Records / Redis Cache section ---------
$ws.Range("B11").Value = 32
$ws.Range("C11").Value = 185
$ws.Range("D11").Value = 45
$ws.Range("E11").Value = 69

$ws.Range("B12").Value = 30
$ws.Range("C12").Value = 184
$ws.Range("D12").Value = 41
$ws.Range("E12").Value = 67

$ws.Range("B14").Value = 57
$ws.Range("C14").Value = 255
$ws.Range("D14").Value = 82
$ws.Range("E14").Value = 108

# --- Updated measurements: Nesting Level / No Cache section ----------------
$ws.Range("B18").Value = 22
$ws.Range("C18").Value = 31
$ws.Range("D18").Value = 27
$ws.Range("E18").Value = 32

$ws.Range("B19").Value = 20
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = 30

$ws.Range("B21").Value = 52
$ws.Range("C21").Value = 69
$ws.Range("D21").Value = 63
$ws.Range("E21").Value = 75

# --- Updated measurements: mount tot data inladen section -------------------
$ws.Range("B26").Value = 41
$ws.Range("C26").Value = 77

# Restore the active selection to where the author left off editing.
$ws.Activate() | Out-Null
$ws.Range("C30").Select() | Out-Null
